# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Appends the 5 new reconciliation rows (95-99) that arrived in this batch,
# extending the used range from A1:I94 to A1:I99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 95 — NOUPA KAMGAING AGNES CHIC MOBILE
$ws.Cells.Item(95, 1).Value = 237675678961
$ws.Cells.Item(95, 2).Value = "NOUPA KAMGAING AGNES CHIC MOBILE"
$ws.Cells.Item(95, 3).Value = 0
$ws.Cells.Item(95, 4).Value = "Ndogbong"
$ws.Cells.Item(95, 5).Value = 10000
$ws.Cells.Item(95, 6).Value = 141
$ws.Cells.Item(95, 7).Value = -9859
$ws.Cells.Item(95, 8).Value = 0.0141
$ws.Cells.Item(95, 9).Value = "Ndogbong"

# Row 96 — Marie Rosine Magne Talla
$ws.Cells.Item(96, 1).Value = 237681678622
$ws.Cells.Item(96, 2).Value = "Marie Rosine Magne Talla"
$ws.Cells.Item(96, 3).Value = "Rte_7"
$ws.Cells.Item(96, 4).Value = "Makepe Conquete"
$ws.Cells.Item(96, 5).Value = 44860
$ws.Cells.Item(96, 6).Value = 49054
$ws.Cells.Item(96, 7).Value = 4194
$ws.Cells.Item(96, 8).Value = 1.093490860454748
$ws.Cells.Item(96, 9).Value = "Ndogbong"

# Row 97 — LA NEGRESSE SARL FOKAM KOM DANICE KEVIN
$ws.Cells.Item(97, 1).Value = 237654041671
$ws.Cells.Item(97, 2).Value = "LA NEGRESSE SARL FOKAM KOM DANICE KEVIN"
$ws.Cells.Item(97, 3).Value = 0
$ws.Cells.Item(97, 4).Value = "Agape Ocm"
$ws.Cells.Item(97, 5).Value = 10000
$ws.Cells.Item(97, 6).Value = 49350
$ws.Cells.Item(97, 7).Value = 39350
$ws.Cells.Item(97, 8).Value = 4.935
$ws.Cells.Item(97, 9).Value = "Cite Sic"

# Row 98 — LA NEGRESSE LTDLA CBOX R0 MALLA TALLA JACQUELINE
$ws.Cells.Item(98, 1).Value = 237675629624
$ws.Cells.Item(98, 2).Value = "LA NEGRESSE LTDLA CBOX R0 MALLA TALLA JACQUELINE"
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = "Esg Building"
$ws.Cells.Item(98, 5).Value = 10000
$ws.Cells.Item(98, 6).Value = 14
$ws.Cells.Item(98, 7).Value = -9986
$ws.Cells.Item(98, 8).Value = 0.0014
$ws.Cells.Item(98, 9).Value = "Ndogbong"

# Row 99 — VAKENA SYLVIE YOK PASL (Pan African Saving and Loan)
$ws.Cells.Item(99, 1).Value = 237683232376
$ws.Cells.Item(99, 2).Value = "VAKENA SYLVIE YOK PASL (Pan African Saving and Loan)"
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = "Bp Cite Ocm"
$ws.Cells.Item(99, 5).Value = 10000
$ws.Cells.Item(99, 6).Value = 500214
$ws.Cells.Item(99, 7).Value = 490214
$ws.Cells.Item(99, 8).Value = 50.0214
$ws.Cells.Item(99, 9).Value = "Cite Sic"
